$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.093.83"
$ws.Range("E2").Value = "  +6.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.528.99"
$ws.Range("E3").Value = "  +3.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "417.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.34%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.657"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.94%  "

# Row 8 - Row 8 -> USDC (was Cardano)
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "

# Row 9 - Row 9 -> Cardano (was USDC)
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.784"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.60%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +20.84%  "

# Row 11 - Avalanche
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.18%  "

# Row 12 - ShibaInu
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000272"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +27.98%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.89%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.088.82"
$ws.Range("E14").Value = "  +2.95%  "

# Row 15 - TRON
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.141"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "

# Row 16 - Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.532.74"
$ws.Range("E17").Value = "  +3.35%  "

# Row 18 - Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "

# Row 19 - Polygon
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.04%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "65.922.37"
$ws.Range("E20").Value = "  +5.88%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "447.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.54%  "

# Row 22 - Litecoin
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "

# Row 23 - ImmutableX
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +1.96%  "

# Row 26 - Filecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.03%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "34.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.23%  "

# Row 28 - Row 28 -> Toncoin (was Cosmos)
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.67%  "

# Row 29 - Row 29 -> Cosmos (was Toncoin)
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.11%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +6.13%  "

# Row 31 - RenderToken
$ws.Range("E31").Value = "  -4.70%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -1.71%  "

# Row 33 - Dai
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.77%  "

# Row 35 - OKB
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.09%  "

# Row 36 - VeChain
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0504"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.63%  "

# Row 37 - PEPE
$ws.Range("D37").Value = "0.0₃0743"
$ws.Range("E37").Value = "  +45.61%  "

# Row 38 - Stellar
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.147"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.05%  "

# Row 39 - FirstDigitalUSD
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "

# Row 40 - Row 40 -> WEMIXToken (was Stacks)
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.80%  "

# Row 41 - Row 41 -> Stacks (was WEMIXToken)
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "

# Row 42 - NEARProtocol
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.01%  "

# Row 43 - Monero
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "

# Row 44 - LidoDAOToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "

# Row 45 - TheGraph
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.311"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.42%  "

# Row 46 - ARBITRUM
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.83%  "

# Row 47 - ThetaToken
$ws.Range("E47").Value = "  -5.98%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +5.88%  "

# Row 49 - Celestia
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.32%  "

# Row 50 - ApeXProtocol
$ws.Range("E50").Value = "  +8.57%  "

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.41%  "
